# Append newly-scraped Lancers listings and bump the "fetched at" timestamp
# on all rows that are still present in this scrape (2026-01-29 12:53:10).
#
# Shape of the change (sheet "ランサーズ" / sheet1):
#   - row 2 (oldest tracked item, 5473940) stays in place, timestamp bumped
#   - two brand-new items (5481801, 5481843) are inserted right after it
#   - the four items that used to be rows 3-6 shift down to rows 5-8,
#     keeping their relative order, each with the refreshed timestamp
#   - one more brand-new item (5481757) is appended as the new last row
#   - column D widens from 28 to 30 characters
#   - every hyperlink in column F is rebuilt so it lines up with its row

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

$timestamp = "2026-01-29 12:53:10"

# --- make room: insert two fresh rows right after row 2 ------------------
$ws.Range("A3:A4").EntireRow.Insert()

# --- drop every existing hyperlink; they get rebuilt from scratch below --
$ws.Hyperlinks.Delete()

# --- widen column D (price) from 28 to 30 characters ----------------------
$ws.Columns.Item(4).ColumnWidth = 29.15

# --- row 2: same listing, just a refreshed fetch timestamp ----------------
$ws.Range("A2").Value = $timestamp

# --- row 3: new listing ----------------------------------------------------
$ws.Range("A3").Value = $timestamp
$ws.Range("B3").Value = "【新規サービス】無人機器と決済連携のシステム開発相談"
$ws.Range("C3").Value = "システム開発"
$ws.Range("D3").Value = "500,000 円 ~ 1,000,000 円 / 固定"
$ws.Range("E3").Value = "期限情報なし"
$ws.Range("F3").Value = "https://www.lancers.jp/work/detail/5481801"
$ws.Range("G3").Value = 125
$ws.Range("H3").Value = "◆開発,システム開発"

# --- row 4: new listing ----------------------------------------------------
$ws.Range("A4").Value = $timestamp
$ws.Range("B4").Value = "【急募】iPadアプリ開発のプロフェッショナルを探しています"
$ws.Range("C4").Value = "システム開発"
$ws.Range("D4").Value = "500,000 円 ~ 1,000,000 円 / 固定"
$ws.Range("E4").Value = "期限情報なし"
$ws.Range("F4").Value = "https://www.lancers.jp/work/detail/5481843"
$ws.Range("G4").Value = 100
$ws.Range("H4").Value = "◆開発 ◇アプリ"

# --- rows 5-8: previously rows 3-6, shifted down, timestamp refreshed -----
$ws.Range("A5").Value = $timestamp
$ws.Range("A6").Value = $timestamp
$ws.Range("A7").Value = $timestamp
$ws.Range("A8").Value = $timestamp

# --- row 9: new listing appended at the end --------------------------------
$ws.Range("A9").Value = $timestamp
$ws.Range("B9").Value = "【急募】BOX内ファイルURLをkintoneに同期するGASプログラム作成"
$ws.Range("C9").Value = "システム開発"
$ws.Range("D9").Value = "5,000 円 ~ 10,000 円 / 固定"
$ws.Range("E9").Value = "期限情報なし"
$ws.Range("F9").Value = "https://www.lancers.jp/work/detail/5481757"
$ws.Range("G9").Value = 10

# --- rebuild every hyperlink in column F, rows 2-9, top to bottom ----------
# (this keeps the relationship ids lined up rId1..rId8 in row order)
$ws.Hyperlinks.Add($ws.Range("F2"), $ws.Range("F2").Value())
$ws.Hyperlinks.Add($ws.Range("F3"), $ws.Range("F3").Value())
$ws.Hyperlinks.Add($ws.Range("F4"), $ws.Range("F4").Value())
$ws.Hyperlinks.Add($ws.Range("F5"), $ws.Range("F5").Value())
$ws.Hyperlinks.Add($ws.Range("F6"), $ws.Range("F6").Value())
$ws.Hyperlinks.Add($ws.Range("F7"), $ws.Range("F7").Value())
$ws.Hyperlinks.Add($ws.Range("F8"), $ws.Range("F8").Value())
$ws.Hyperlinks.Add($ws.Range("F9"), $ws.Range("F9").Value())

$ws.Range("F2:F9").Style = "Hyperlink"
